$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) cells whose new values look numeric,
# so Excel stores them as text (matching the original inlineStr type)
# instead of auto-converting to a number.
$ws.Range('D2').Value = '36.609.78'
$ws.Range('E2').Value = '  -0.21%  '
$ws.Range('D3').Value = '1.972.53'
$ws.Range('E3').Value = '  +0.34%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.19'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.00%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.628'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.83%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '59.96'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.50%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.378'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.69%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0788'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.03%  '
$ws.Range('E11').Value = '  +0.40%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.21'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.99%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.839'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.97%  '
$ws.Range('D14').Value = '2.267.70'
$ws.Range('E14').Value = '  +0.48%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.71'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.96%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.42'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.79%  '
$ws.Range('D17').Value = '1.972.00'
$ws.Range('E17').Value = '  +0.58%  '
$ws.Range('D18').Value = '36.564.58'
$ws.Range('E18').Value = '  -0.14%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '69.77'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.11%  '
$ws.Range('D20').Value = '0.0₃0854'
$ws.Range('E20').Value = '  -0.81%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.10'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.12%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '229.52'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.46%  '
$ws.Range('E23').Value = '  +0.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.44'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.23%  '
$ws.Range('E25').Value = '  +1.45%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.146'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +6.19%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.16'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.66%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '161.54'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.74%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.36'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.12%  '
$ws.Range('E30').Value = '  +19.39%  '
$ws.Range('E31').Value = '  +1.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.79'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.77%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0616'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.41%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.50'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.83%  '
$ws.Range('B35').Value = 'BinanceUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.17%  '
$ws.Range('B36').Value = 'LidoDAOToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.27'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.75%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.77'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.49%  '
$ws.Range('E38').Value = '  -3.50%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.40'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -11.33%  '
$ws.Range('E40').Value = '  -2.37%  '
$ws.Range('E41').Value = '  +0.53%  '
$ws.Range('E42').Value = '  +0.08%  '
$ws.Range('E43').Value = '  -0.93%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '15.84'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.03%  '
$ws.Range('D45').Value = '1.364.01'
$ws.Range('E45').Value = '  -0.53%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '88.97'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.78%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.02'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.51%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.19'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.65%  '
$ws.Range('E49').Value = '  -0.68%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '46.04'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +5.96%  '
$ws.Range('D51').Value = '2.161.51'
$ws.Range('E51').Value = '  +0.63%  '
